$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Simple value edits
$ws.Range("C13").Value = 340926355.5
$ws.Range("C14").Value = 985666.55
$ws.Range("C15").Value = 300000000
$ws.Range("C16").Value = -51786608.799999997
$ws.Range("C19").Value = -392700000
$ws.Range("C22").Value = -33376095

# Convert C18 and C21 from literal values to formulas (mirroring the other columns' pattern)
$ws.Range("C18").Formula = "=SUM(C12:C17)"
$ws.Range("C21").Formula = "=SUM(C18:C20)"

$wb.Application.Calculate()
